{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, etc.)\n// in bold + a dark slate color (#2C3E50) across specific resume bullet\n// paragraphs, matching the \"hybrid bold + color highlighting\" commit.\n//\n// Strategy: for each target paragraph (identified by a unique, stable text\n// prefix), locate each metric substring with Paragraph.search() (which\n// returns a Range scoped to just that paragraph) and apply bold + color to\n// that Range's font. Word/Office.js automatically splits the run(s) around\n// the matched sub-range, which reproduces the same run layout shown in the\n// diff (plain-text runs around freshly bolded/colored metric runs).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Each entry: [unique paragraph-start prefix, ordered list of metric\n// substrings to bold+color within that paragraph]. Order matters when a\n// paragraph contains more than one highlighted metric so they come out in\n// left-to-right reading order, matching the diff.\nconst plan = [\n  [\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters\",\n    [\"23%\", \"64%\"],\n  ],\n  [\n    \"\u2022 Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2%\",\n    [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"],\n  ],\n  [\n    \"\u2022 Trigonometric algorithm for boundary estimation\",\n    [\"73.5%\", \"$4.7M\"],\n  ],\n  [\n    \"\u2022 Built real-time FEC analysis systems\",\n    [\"$2\"],\n  ],\n  [\n    \"\u2022 Modernized legacy ETL processes\",\n    [\"57%\"],\n  ],\n  [\n    \"\u2022 Predictive excellence: Utilized advanced sampling methods\",\n    [\"\u00b14.2%\", \"\u00b12.1%\"],\n  ],\n  [\n    \"\u2022 Increased voter turnout prediction accuracy from 71%\",\n    [\"71%\", \"87%\"],\n  ],\n  [\n    \"\u2022 Methodological advancement\",\n    [\"34%\", \"28%\"],\n  ],\n];\n\nconst items = paragraphs.items;\n\nfor (const [prefix, terms] of plan) {\n  const para = items.find((p) => p.text.startsWith(prefix));\n  if (!para) continue;\n\n  for (const term of terms) {\n    const found = para.search(term, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n\n    // Bold + color every occurrence of the metric inside this paragraph\n    // (normally exactly one), leaving the surrounding text untouched.\n    found.items.forEach((hit) => {\n      hit.font.bold = true;\n      hit.font.color = \"#2C3E50\";\n    });\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, etc.)\n# in bold + a dark slate color (#2C3E50) across specific resume bullet\n# paragraphs, matching the \"hybrid bold + color highlighting\" commit.\n#\n# Strategy: for each target paragraph (identified by a unique, stable text\n# prefix), re-fetch a fresh Range over that whole paragraph for every metric\n# and use Range.Find.Execute to collapse it onto just that metric's text,\n# then set Font.Bold / Font.Color on the collapsed range. Word automatically\n# splits the surrounding run(s) around the newly formatted sub-range, which\n# reproduces the same run layout shown in the diff (plain-text runs around\n# freshly bolded/colored metric runs).\n\n$d = $word.ActiveDocument\n\n# Each entry: paragraph-start prefix -> ordered list of metric substrings to\n# bold+color within that paragraph. Order matters when a paragraph contains\n# more than one highlighted metric so they come out in left-to-right reading\n# order, matching the diff.\n$plan = @(\n    @{ Prefix = \"\u2022 Discovered systematic race coding errors\"; Terms = @(\"23%\", \"64%\") },\n    @{ Prefix = \"\u2022 Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2%\"; Terms = @(\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\") },\n    @{ Prefix = \"\u2022 Trigonometric algorithm for boundary estimation\"; Terms = @(\"73.5%\", \"$4.7M\") },\n    @{ Prefix = \"\u2022 Built real-time FEC analysis systems\"; Terms = @(\"$2\") },\n    @{ Prefix = \"\u2022 Modernized legacy ETL processes\"; Terms = @(\"57%\") },\n    @{ Prefix = \"\u2022 Predictive excellence: Utilized advanced sampling methods\"; Terms = @(\"\u00b14.2%\", \"\u00b12.1%\") },\n    @{ Prefix = \"\u2022 Increased voter turnout prediction accuracy from 71%\"; Terms = @(\"71%\", \"87%\") },\n    @{ Prefix = \"\u2022 Methodological advancement\"; Terms = @(\"34%\", \"28%\") }\n)\n\nforeach ($entry in $plan) {\n    $prefix = $entry.Prefix\n    $terms = $entry.Terms\n\n    $target = $null\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Text.StartsWith($prefix)) {\n            $target = $p\n            break\n        }\n    }\n    if ($null -eq $target) {\n        continue\n    }\n\n    foreach ($t in $terms) {\n        $rng = $target.Range\n        $rng.Find.ClearFormatting()\n        $rng.Find.MatchCase = $true\n        $rng.Find.MatchWholeWord = $false\n        $rng.Find.MatchWildcards = $false\n        $found = $rng.Find.Execute($t)\n        if ($found) {\n            $rng.Font.Bold = 1\n            $rng.Font.Color = \"2C3E50\"\n        }\n    }\n}\n"}
